# Apply the "cadastro" update: change the Usuario value on the Cadastro sheet
# and move the active selection to A3, as reflected in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cadastro")

# Update the username value in A2 (was "pradov1027")
$ws.Range("A2").Value = "pradov1033"

# Activate the Cadastro sheet and move the selection to A3
$ws.Activate()
$ws.Range("A3").Select()
